$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new log entry was recorded at the top of the "Loan" (September) block.
# Insert a blank row above the existing row 27, pushing the current rows
# 27-35 down to 28-36 (Excel shifts formulas/data automatically), then
# populate the newly inserted row 27 with the new log entry.
$ws.Rows("27:27").Insert()

$ws.Range("R27").Value = "broker"
$ws.Range("S27").Value = "2024-09-01 22:35:38"
